$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FamiliarStats_Master")

# --- Restructure columns on FamiliarStats_Master ---
# Old layout:  A RowName | B SummonCost | C AttackTypeTag | D RoleTypeTag | E AttackInterval |
#              F AttackRange | G BaseMaxHP | H BaseMaxMP | I BaseAttackPower | J BaseDefense |
#              K BaseMoveSpeed | L BaseCritRate
# New layout:  A RowName | B SummonCost | C AttackTypeTag | D RoleTypeTag | E AttackSpeed |
#              F AttackRange | G Cooldown | H BaseMaxHP | I BaseMaxMP | J BaseAttackPower |
#              K BaseDefense | L BaseMoveSpeed | M BaseCritRate

# Drop the old "AttackInterval" column entirely (column E).
$ws.Columns("E:E").Delete()

# Insert a fresh column E for the new "AttackSpeed" stat.
$ws.Columns("E:E").Insert()

# Insert a fresh column G for the new "Cooldown" stat (after AttackRange, which stays in F).
$ws.Columns("G:G").Insert()

# --- Headers (row 1) ---
$ws.Range("E1").Value2 = "AttackSpeed"
$ws.Range("G1").Value2 = "Cooldown"

# --- Sample data (row 2) ---
$ws.Range("E2").Value2 = 0.1
$ws.Range("G2").Value2 = 0.2

# Column F ("AttackRange") picked up an explicit width in the authored workbook.
$ws.Columns("F:F").ColumnWidth = 8.43

# --- Selection matches the post-edit workbook state ---
$ws.Range("G2").Select()

# --- Data validation for the two new columns ---
$rngCooldown = $ws.Range("G2:G1048576")
$rngCooldown.Validation.Add(2, 1, 7, 0)
$rngCooldown.Validation.ErrorTitle = "단위 확인 (Cooldown)"
$rngCooldown.Validation.ErrorMessage = "쿨타임은 0보다 작을 수 없습니다._x000a__x000a_※ 단위: 초 (Seconds)_x000a_- 1분 = '60' 입력_x000a_- 0.5초 = '0.5' 입력_x000a_- 쿨타임 없음 = '0' 입력_x000a__x000a_밀리초(ms) 단위가 아닙니다! 다시 확인해주세요."
$rngCooldown.Validation.ShowInput = $true
$rngCooldown.Validation.ShowError = $true

$rngAttackSpeed = $ws.Range("E2:E1048576")
$rngAttackSpeed.Validation.Add(2, 1, 5, 0)
$rngAttackSpeed.Validation.ErrorTitle = "입력 오류"
$rngAttackSpeed.Validation.ErrorMessage = "음수(-)와 0은 입력할 수 없습니다. 0 보다 큰 숫자를 입력해주세요."
$rngAttackSpeed.Validation.ShowInput = $true
$rngAttackSpeed.Validation.ShowError = $true

# --- Drop the stale external workbook link (EnemyStats_Master.xlsx) ---
$wb.BreakLink("EnemyStats_Master.xlsx", 1)
